$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row that contains "Change" (row 80) so that the
# remaining skill list shifts up by one row, and Excel drops the now
# unused "Change" entry from the shared strings table automatically.
$ws.Rows.Item(80).Delete()

# Restore a sensible selection/view similar to what Excel leaves behind
# after deleting a row via the row header context menu.
$ws.Application.Goto($ws.Range("A80:XFD80"), $false)
$ws.Application.ActiveWindow.ScrollRow = 179
